$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3016.0667
$ws.Range("I33").Value = 3302.5833
$ws.Range("K33").Value = 3302.5833
$ws.Range("M33").Value = -3073.5833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 57492.5
$ws.Range("J105").Value = 57492.5
$ws.Range("L105").Value = 57492.5
$ws.Range("N105").Value = -64480.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4284.9546
$ws.Range("I132").Value = 4128.75
$ws.Range("K132").Value = 12386.25
$ws.Range("M132").Value = -9856.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3327
$ws.Range("I138").Value = 2257.6667
$ws.Range("K138").Value = 6773.000100000001
$ws.Range("M138").Value = -1633.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1207.0769
$ws.Range("I2").Value = 1219.3
$ws.Range("J2").Value = 1166.3334
$ws.Range("K2").Value = 1219.3
$ws.Range("L2").Value = 1166.3334
$ws.Range("M2").Value = -1106.3
$ws.Range("N2").Value = -1392.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 55556956
$ws.Range("I45").Value = 71429800
$ws.Range("K45").Value = 71429800
$ws.Range("M45").Value = -71429423

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1207.0769
$ws.Range("I116").Value = 1219.3
$ws.Range("J116").Value = 1166.3334
$ws.Range("K116").Value = 1219.3
$ws.Range("L116").Value = 1166.3334
$ws.Range("M116").Value = 1074.7
$ws.Range("N116").Value = -5754.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1207.0769
$ws.Range("I3").Value = 1219.3
$ws.Range("J3").Value = 1166.3334
$ws.Range("K3").Value = 1219.3
$ws.Range("L3").Value = 1166.3334
$ws.Range("M3").Value = -1105.3
$ws.Range("N3").Value = -1394.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2153.8215
$ws.Range("I86").Value = 1532.3529
$ws.Range("J86").Value = 3114.2727
$ws.Range("K86").Value = 1532.3529
$ws.Range("L86").Value = 3114.2727
$ws.Range("M86").Value = -409.3529000000001
$ws.Range("N86").Value = -5360.2727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2153.8215
$ws.Range("I89").Value = 1532.3529
$ws.Range("J89").Value = 3114.2727
$ws.Range("K89").Value = 7661.7645
$ws.Range("L89").Value = 15571.3635
$ws.Range("M89").Value = -2045.7645
$ws.Range("N89").Value = -26803.3635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9526.3125
$ws.Range("I99").Value = 21942
$ws.Range("K99").Value = 21942
$ws.Range("M99").Value = -20444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 627135.4
$ws.Range("I134").Value = 1025.9231
$ws.Range("J134").Value = 3340276.2
$ws.Range("K134").Value = 3077.7693
$ws.Range("L134").Value = 10020828.6
$ws.Range("M134").Value = -542.7692999999999
$ws.Range("N134").Value = -10025898.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 7259.4
$ws.Range("J28").Value = 7259.4
$ws.Range("L28").Value = 7259.4
$ws.Range("N28").Value = -7749.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 639599.9
$ws.Range("I31").Value = 11236
$ws.Range("K31").Value = 11236
$ws.Range("M31").Value = -10941

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 639599.9
$ws.Range("I34").Value = 11236
$ws.Range("K34").Value = 11236
$ws.Range("M34").Value = -11034

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1113168.8
$ws.Range("I134").Value = 1252037.1
$ws.Range("K134").Value = 3756111.3
$ws.Range("M134").Value = -3753576.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6100056.5
$ws.Range("I4").Value = 7117705.5
$ws.Range("K4").Value = 21353116.5
$ws.Range("M4").Value = -21353004.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3297.8823
$ws.Range("I80").Value = 2865
$ws.Range("J80").Value = 3534
$ws.Range("K80").Value = 8595
$ws.Range("L80").Value = 10602
$ws.Range("M80").Value = -7659
$ws.Range("N80").Value = -12474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3297.8823
$ws.Range("I83").Value = 2865
$ws.Range("J83").Value = 3534
$ws.Range("K83").Value = 25785
$ws.Range("L83").Value = 31806
$ws.Range("M83").Value = -21105
$ws.Range("N83").Value = -41166

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1970.3572
$ws.Range("J122").Value = 2064.2307
$ws.Range("L122").Value = 18578.0763
$ws.Range("N122").Value = -23478.0763

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 447951
$ws.Range("I128").Value = 447951
$ws.Range("K128").Value = 1343853
$ws.Range("M128").Value = -1338873

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5082.8125
$ws.Range("I131").Value = 600
$ws.Range("J131").Value = 5723.2144
$ws.Range("K131").Value = 1800
$ws.Range("L131").Value = 17169.6432
$ws.Range("M131").Value = 3240
$ws.Range("N131").Value = -27249.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2237.524
$ws.Range("I132").Value = 1998.5
$ws.Range("J132").Value = 2384.6155
$ws.Range("K132").Value = 17986.5
$ws.Range("L132").Value = 21461.5395
$ws.Range("M132").Value = -15456.5
$ws.Range("N132").Value = -26521.5395

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 50001264
$ws.Range("I133").Value = 50001264
$ws.Range("K133").Value = 150003792
$ws.Range("M133").Value = -149998732

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3537722.8
$ws.Range("J11").Value = 5893715.5
$ws.Range("L11").Value = 5893715.5
$ws.Range("N11").Value = -5893993.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.5
$ws.Range("I80").Value = 2499
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2499
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1501
$ws.Range("N80").Value = -4996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2749.5
$ws.Range("I83").Value = 2499
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 12495
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -7503
$ws.Range("N83").Value = -24984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 58834212
$ws.Range("I132").Value = 90912140
$ws.Range("J132").Value = 24668.666
$ws.Range("K132").Value = 272736420
$ws.Range("L132").Value = 74005.99800000001
$ws.Range("M132").Value = -272733890
$ws.Range("N132").Value = -79065.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 887246.5600000001
$ws.Range("I43").Value = 45945
$ws.Range("K43").Value = 45945
$ws.Range("M43").Value = -45752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 19994.5
$ws.Range("J76").Value = 19994.5
$ws.Range("L76").Value = 19994.5
$ws.Range("N76").Value = -20670.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 19994.5
$ws.Range("J79").Value = 19994.5
$ws.Range("L79").Value = 19994.5
$ws.Range("N79").Value = -22334.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 15000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 20666.666
$ws.Range("I29").Value = 3500
$ws.Range("K29").Value = 3500
$ws.Range("M29").Value = -3210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5063.273
$ws.Range("I122").Value = 4966.222
$ws.Range("K122").Value = 14898.666
$ws.Range("M122").Value = -12448.666
